# Finalized market cap csv
# Swap the "Market cap" / "Change" column headers: column B now reads
# "MarketCap" and column C (the percentage-formatted column) now reads
# "Change".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "MarketCap"
$ws.Range("C1").Value = "Change"

# Give column B an explicit width (author widened it to fit "MarketCap").
$ws.Columns.Item(2).ColumnWidth = 12.6

# Leave the selection on E11, matching the author's last on-screen selection.
$ws.Range("E11").Select() | Out-Null
